$d = $word.ActiveDocument

# Find the target paragraph (end of the bullet list item about random generation ranges)
$marker = "I may not get the same answer due to the number of random generation ranges. And also row shuffling (e.g. the train_test_split function is shuffle true by default)."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*$marker*") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after the target paragraph, collapsed to its end
$insertRange = $target.Range
$insertRange.Collapse(0)  # wdCollapseEnd
$insertRange.InsertParagraphAfter()

# Move to the newly created paragraph and set its text + list formatting
$newPara = $target.Next()
$newPara.Range.Text = "There will be unused code. In production level quality, there should not be any unused code."

# Match the numbered/bulleted list formatting (numId 3) of the preceding list items
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($target.Range.ListFormat.ListTemplate)
